# Insert a new price-report row at row 40 (pushing the existing rows 40-66
# down to 41-67) and populate it with the new "Superior Seedless" record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 40..66 down one position, creating a blank row 40.
$ws.Rows.Item(40).Insert()

# Fill in the new row 40 with the reported values.
$ws.Range("A40").Value = 8
$ws.Range("B40").Value = "Terminal La Palmera de La Serena"
$ws.Range("C40").Value = "Coquimbo"
$ws.Range("D40").Value = 44582
$ws.Range("E40").Value = 4
$ws.Range("F40").Value = "Fruta"
$ws.Range("G40").Value = 100109
$ws.Range("H40").Value = "Uva"
$ws.Range("I40").Value = 100109001
$ws.Range("J40").Value = "Uva"
$ws.Range("K40").Value = "Superior Seedless"
$ws.Range("L40").Value = "Primera"
$ws.Range("M40").Value = 400
$ws.Range("N40").Value = 12000
$ws.Range("O40").Value = 13000
$ws.Range("P40").Value = 12500
$ws.Range("Q40").Value = "`$/bandeja 18 kilos"
$ws.Range("R40").Value = "Provincia del Elquí"
$ws.Range("S40").Value = 694
$ws.Range("T40").Value = 18
